$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-13 00:57:49"
$wsZh.Range("H2").Value = "2016-03-13 00:58:06"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-13 00:57:56"
$wsDe.Range("H2").Value = "2016-03-13 00:58:12"
